$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# New row 111
Set-TextCell 111 1 "2025-10-26 02:30:05"
Set-TextCell 111 2 "Noah Dubitzky"
$ws.Cells.Item(111, 3).Value = 8450689526
Set-TextCell 111 4 "13052054965"
Set-TextCell 111 5 "10:30"
Set-TextCell 111 6 ""
Set-TextCell 111 7 ""

# New row 112
Set-TextCell 112 1 "2025-10-25 22:45:50"
Set-TextCell 112 2 "Noah Dubitzky"
$ws.Cells.Item(112, 3).Value = 8450689526
Set-TextCell 112 4 "13052054965"
Set-TextCell 112 5 "10:45"
Set-TextCell 112 6 ""
Set-TextCell 112 7 ""
